$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 315
$ws1.Range("F4").Value = 8265
$ws1.Range("F5").Value = 6027
$ws1.Range("F6").Value = 516
$ws1.Range("F7").Value = 101
$ws1.Range("F11").Value = 935
$ws1.Range("F12").Value = 80

# Sheet "全部类型" (all types) - same rows updated with the same new values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 315
$ws4.Range("F4").Value = 8265
$ws4.Range("F5").Value = 6027
$ws4.Range("F6").Value = 516
$ws4.Range("F7").Value = 101
$ws4.Range("F15").Value = 935
$ws4.Range("F16").Value = 80

$wb.Save()
